$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.510.52"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.915.78"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.53"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4819"
$ws.Range("E7").Value = "  +2.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2890"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06738"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "111.17"
$ws.Range("E10").Value = "  +3.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.27"
$ws.Range("E11").Value = "  +4.85%  "
$ws.Range("D12").Value = "1.919.21"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.260"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6717"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "288.04"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "30.532.52"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007623"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9992"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.91"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "2.167.24"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.497"
$ws.Range("E22").Value = "  +4.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.424"
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.472"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.30"
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.37"
$ws.Range("E27").Value = "  -5.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.145"
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1062"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.409"
$ws.Range("E30").Value = "  +3.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.169"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.039"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04995"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7294"
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.134"
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02056"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9989"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.720"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.668"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.57"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.013"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4436"
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8647"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.887"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "68.03"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.361"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.78"
$ws.Range("E48").Value = "  -3.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.337"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1242"
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.88"
$ws.Range("E51").Value = "  -0.21%  "
